$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.424.67'
$ws.Range("E2").Value = '  -3.78%  '

$ws.Range("D3").Value = '2.618.00'
$ws.Range("E3").Value = '  -3.43%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '522.03'
$ws.Range("E5").Value = '  -1.21%  '

$ws.Range("D6").Value = '142.56'
$ws.Range("E6").Value = '  -3.05%  '

$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("E8").Value = '  -1.78%  '

$ws.Range("E9").Value = '  -7.08%  '

$ws.Range("D10").Value = '0.103'
$ws.Range("E10").Value = '  -2.57%  '

$ws.Range("D11").Value = '0.334'
$ws.Range("E11").Value = '  -1.77%  '

$ws.Range("E12").Value = '  +0.89%  '

$ws.Range("D13").Value = '3.078.24'
$ws.Range("E13").Value = '  -3.46%  '

$ws.Range("D14").Value = '58.397.54'
$ws.Range("E14").Value = '  -3.83%  '

$ws.Range("D15").Value = '20.95'
$ws.Range("E15").Value = '  -2.19%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0000136'
$ws.Range("E16").Value = '  -1.85%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.632.86'
$ws.Range("E17").Value = '  -3.62%  '

$ws.Range("D18").Value = '337.49'
$ws.Range("E18").Value = '  -2.02%  '

$ws.Range("E19").Value = '  -2.21%  '

$ws.Range("D20").Value = '10.38'
$ws.Range("E20").Value = '  -1.68%  '

$ws.Range("E21").Value = '  -2.50%  '

$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").Value = '65.05'
$ws.Range("E23").Value = '  +2.65%  '

$ws.Range("E24").Value = '  -0.95%  '

$ws.Range("E25").Value = '  -2.58%  '

$ws.Range("E26").Value = '  +0.35%  '

$ws.Range("E27").Value = '  -2.55%  '

$ws.Range("D28").Value = '0.0₃0790'
$ws.Range("E28").Value = '  -3.99%  '

$ws.Range("E29").Value = '  -3.19%  '

$ws.Range("E31").Value = '  -0.72%  '

$ws.Range("D32").Value = '18.75'
$ws.Range("E32").Value = '  -1.67%  '

$ws.Range("D33").Value = '150.20'
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("E34").Value = '  -3.60%  '

$ws.Range("E35").Value = '  -3.21%  '

$ws.Range("D36").Value = '0.895'
$ws.Range("E36").Value = '  -2.99%  '

$ws.Range("D37").Value = '0.850'
$ws.Range("E37").Value = '  -5.68%  '

$ws.Range("E38").Value = '  -2.84%  '

$ws.Range("E39").Value = '  -6.48%  '

$ws.Range("E40").Value = '  -1.55%  '

$ws.Range("E41").Value = '  +0.36%  '

$ws.Range("E42").Value = '  -3.60%  '

$ws.Range("D43").Value = '0.0971'
$ws.Range("E43").Value = '  -1.57%  '

$ws.Range("D44").Value = '268.10'
$ws.Range("E44").Value = '  -4.49%  '

$ws.Range("D46").Value = '19.06'
$ws.Range("E46").Value = '  -5.65%  '

$ws.Range("D47").Value = '0.0531'
$ws.Range("E47").Value = '  -2.07%  '

$ws.Range("D48").Value = '2.028.03'
$ws.Range("E48").Value = '  -3.79%  '

$ws.Range("D49").Value = '0.0229'
$ws.Range("E49").Value = '  -1.35%  '

$ws.Range("E50").Value = '  -7.69%  '

$ws.Range("E51").Value = '  -5.28%  '
